$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('B1').Value = 'Honda HR-V'
$ws.Range('C1').Value = 'Dacia Jogger'
$ws.Range('D1').Value = 'DS 4'
$ws.Range('E1').Value = 'Volvo C40 Recharge'
$ws.Range('F1').Value = 'Mercedes-Benz C-Class'
$ws.Range('G1').Value = 'Kia EV6'
$ws.Range('H1').Value = 'VW Multivan'
$ws.Range('I1').Value = 'Peugeot 308'
$ws.Range('J1').Value = 'ORA FUNKY CAT'
$ws.Range('K1').Value = 'WEY Coffee 01'
$ws.Range('L1').Value = 'Kia Niro'
$ws.Range('M1').Value = 'Genesis GV60'
$ws.Range('N1').Value = 'Tesla Model Y'
$ws.Range('O1').Value = 'Land Rover Discovery Sport'
$ws.Range('P1').Value = 'CHERY OMODA5'
$ws.Range('Q1').Value = 'Volkswagen ID. Buzz'
$ws.Range('R1').Value = 'Volkswagen Touran'
$ws.Range('S1').Value = 'Lexus RX'
$ws.Range('T1').Value = 'Skoda Octavia'
$ws.Range('U1').Value = 'MG 4 Electric'
$ws.Range('V1').Value = 'Ford Ranger'
$ws.Range('W1').Value = 'Lucid Air'
$ws.Range('X1').Value = 'Jeep Grand Cherokee'
$ws.Range('Y1').Value = 'Volkswagen Amarok'
$ws.Range('Z1').Value = 'Mercedes-Benz GLC'
$ws.Range('AA1').Value = 'Maxus MIFA 9'
$ws.Range('AB1').Value = 'Ford Puma'
$ws.Range('AC1').Value = 'BMW 2 Series Coupé'
$ws.Range('AD1').Value = 'Renault Megane E-Tech'
$ws.Range('AE1').Value = 'Volkswagen Polo'
$ws.Range('AF1').Value = 'Lexus NX'
$ws.Range('AG1').Value = 'Volkswagen Taigo'
$ws.Range('AH1').Value = 'Nissan Ariya'
$ws.Range('AI1').Value = 'smart #1'
$ws.Range('AJ1').Value = 'Hyundai IONIQ 6'
$ws.Range('AK1').Value = 'Toyota Corolla Cross'
$ws.Range('AL1').Value = 'Land Rover Range Rover Sport'
$ws.Range('AM1').Value = 'Isuzu D-Max Crew Cab'
$ws.Range('AN1').Value = 'NIO ET7'
$ws.Range('AO1').Value = 'Land Rover Range Rover'
$ws.Range('AP1').Value = 'Renault Austral'
$ws.Range('AQ1').Value = 'DS 9'
$ws.Range('AR1').Value = 'Tesla Model S'
$ws.Range('AS1').Value = 'Honda Civic'
$ws.Range('AT1').Value = 'Nissan X-Trail'
$ws.Range('AU1').Value = 'WEY Coffee 02'
$ws.Range('AV1').Value = 'Toyota bZ4X'
$ws.Range('AW1').Value = 'Seat Ibiza'
$ws.Range('AX1').Value = 'BMW X1'
$ws.Range('AY1').Value = 'Mobilize Limo'
$ws.Range('AZ1').Value = 'Mercedes-EQ EQE'
$ws.Range('BA1').Value = 'BYD Atto 3'
$ws.Range('BB1').Value = 'Citroën C5 X'
$ws.Range('BC1').Value = 'SEAT Arona'
$ws.Range('BD1').Value = 'Mazda MAZDA CX-60'
$ws.Range('BE1').Value = 'BMW 2 Series Active Tourer'
$ws.Range('BF1').Value = 'Volkswagen Golf'
$ws.Range('BG1').Value = 'Kia Sportage'
$ws.Range('BH1').Value = 'BMW i4'
$ws.Range('BI1').Value = 'Mercedes-Benz T-Class'
$ws.Range('BJ1').Value = 'Toyota Aygo X'
$ws.Range('BK1').Value = 'Alfa Romeo Tonale'
$ws.Range('BL1').Value = 'Cupra Born'
